# Update crypto "Price" (column D) and "Volume(1h)" (column E) values to
# the freshly-scraped figures from the Sat Aug 3 10:51:17 UTC 2024 run.
#
# Column E values must keep their two-space padding on both sides,
# matching the original formatting (e.g. "  -4.16%  ").
#
# Column D values are plain text in the source data (they come from a
# scraper, not user typing), and some of them look like ordinary decimal
# numbers (e.g. "540.70", "12.05"). Excel's normal cell-entry behavior
# would silently convert such text into a numeric value and drop the
# formatting (trailing zero, etc.), so for any new D value that would be
# auto-recognized as a plain number we force the cell to Text format
# first so the exact original string is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new D value (or $null to leave unchanged), new E percent text (or $null to leave unchanged)
$updates = @(
    @{ Row = 2;  D = "61.683.50"; E = "-4.16%" },
    @{ Row = 3;  D = "2.979.92";  E = "-5.07%" },
    @{ Row = 4;  D = $null;       E = "+0.01%" },
    @{ Row = 5;  D = "540.70";    E = "-5.44%" },
    @{ Row = 6;  D = "151.65";    E = "-7.64%" },
    @{ Row = 8;  D = $null;       E = "-1.33%" },
    @{ Row = 9;  D = "2.991.48";  E = "-5.10%" },
    @{ Row = 10; D = $null;       E = "-3.67%" },
    @{ Row = 11; D = $null;       E = "-7.13%" },
    @{ Row = 12; D = "0.369";     E = "-3.91%" },
    @{ Row = 13; D = "3.501.59";  E = "-5.04%" },
    @{ Row = 14; D = $null;       E = "-2.27%" },
    @{ Row = 15; D = "61.732.92"; E = "-4.08%" },
    @{ Row = 16; D = "23.86";     E = "-4.54%" },
    @{ Row = 17; D = "2.979.93";  E = "-5.16%" },
    @{ Row = 18; D = $null;       E = "-5.69%" },
    @{ Row = 19; D = $null;       E = "-1.82%" },
    @{ Row = 20; D = "12.05";     E = "-3.96%" },
    @{ Row = 21; D = "381.55";    E = "-5.91%" },
    @{ Row = 22; D = $null;       E = "-5.37%" },
    @{ Row = 23; D = $null;       E = "+0.05%" },
    @{ Row = 24; D = "5.65";      E = "-3.55%" },
    @{ Row = 25; D = "65.91";     E = "-4.36%" },
    @{ Row = 26; D = $null;       E = "-2.55%" },
    @{ Row = 27; D = "3.102.95";  E = "-5.25%" },
    @{ Row = 28; D = "0.190";     E = "-2.46%" },
    @{ Row = 29; D = "0.998";     E = "+0.26%" },
    @{ Row = 30; D = "0.0₃0939";  E = "-7.87%" },
    @{ Row = 31; D = "8.17";      E = "-7.97%" },
    @{ Row = 32; D = $null;       E = "+0.02%" },
    @{ Row = 33; D = $null;       E = "-4.39%" },
    @{ Row = 34; D = "20.49";     E = "-3.48%" },
    @{ Row = 35; D = "159.84";    E = "-1.84%" },
    @{ Row = 36; D = "4.58";      E = "-5.81%" },
    @{ Row = 37; D = $null;       E = "-6.05%" },
    @{ Row = 38; D = $null;       E = "-5.03%" },
    @{ Row = 39; D = $null;       E = "-6.48%" },
    @{ Row = 40; D = $null;       E = "-8.32%" },
    @{ Row = 41; D = "37.57";     E = $null },
    @{ Row = 42; D = "2.420.29";  E = "-8.15%" },
    @{ Row = 43; D = $null;       E = "-4.61%" },
    @{ Row = 44; D = "22.03";     E = "-6.85%" },
    @{ Row = 45; D = "0.671";     E = "-2.80%" },
    @{ Row = 46; D = "0.0590";    E = "-3.48%" },
    @{ Row = 47; D = "5.12";      E = "-4.98%" },
    @{ Row = 48; D = $null;       E = "+0.08%" },
    @{ Row = 49; D = $null;       E = "-3.81%" },
    @{ Row = 50; D = $null;       E = "-2.29%" },
    @{ Row = 51; D = $null;       E = "-6.91%" }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        $dRange = $ws.Range("D$r")
        # If the new text would otherwise be auto-parsed as a plain number
        # by Excel (losing its exact textual form), pin the cell to Text
        # format first so the literal string is stored unchanged.
        if ($u.D -match '^[0-9]+(\.[0-9]+)?$') {
            $dRange.NumberFormat = "@"
        }
        $dRange.Value = $u.D
    }

    if ($null -ne $u.E) {
        $ws.Range("E$r").Value = "  " + $u.E + "  "
    }
}
